# Update SAM ToDo and Beta feedback spreadsheets
# Adds a new feedback row (row 56) to Sheet1 of the SAM 2014 Beta feedback log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from the row above (row 55) into the new row 56 so the
# existing shared styles (date format, wrap-text, etc.) are reused instead
# of creating new style entries.
$ws.Range("A55").Copy()
$ws.Range("A56").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B55").Copy()
$ws.Range("B56").PasteSpecial(-4122)
$ws.Range("C55").Copy()
$ws.Range("C56").PasteSpecial(-4122)
$ws.Range("D55").Copy()
$ws.Range("D56").PasteSpecial(-4122)
$ws.Range("F55").Copy()
$ws.Range("F56").PasteSpecial(-4122)

# New feedback entry data
$ws.Range("A56").Value = 41934
$ws.Range("B56").Value = "Forum"
$ws.Range("C56").Value = "Stephen.Frank@nrel.gov"
$ws.Range("D56").Value = "Request to make it easier to transfer sam inputs to SDK -- enhancements to inputs browser"
$ws.Range("F56").Value = 41934

# Match the row height used for the new (wrapped, 2-line) entry
$ws.Rows.Item(56).RowHeight = 30

# Move the selection to the next empty row beneath the new entry, and
# clear the application's horizontal scroll position of the top sheet view
# by selecting in column A.
$ws.Range("A57").Select() | Out-Null
